$wb = $excel.ActiveWorkbook

# OFF sheet - Road row (row 3) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 145
$wsOff.Range("C3").Value = 95
$wsOff.Range("D3").Value = 28
$wsOff.Range("E3").Value = 17

# DEF sheet - Road row (row 3) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 152
$wsDef.Range("C3").Value = 101
$wsDef.Range("D3").Value = 33
$wsDef.Range("F3").Value = 4
